# Generate Report for Handoff
#
# The e895b03e-a365-48c4-96fa-0f40707d4535.md file just had a new handoff
# generated, so its "Latest Handoff" timestamps move forward on every sheet
# that tracks it:
#   - Overview!D5            ("Latest Handoff Date")       -> 2016-03-22 05:04:12
#   - zh-cn!E5  (row 5 data) ("Latest Handoff Datetime")   -> 2016-03-22 05:04:05
#   - de-de!E5  (row 5 data) ("Latest Handoff Datetime")   -> 2016-03-22 05:04:12

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-03-22 05:04:12"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-22 05:04:05"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-22 05:04:12"
